$d = $word.ActiveDocument

# Locate the paragraph containing the standalone "Test." sentence and
# remove just its text (keep the paragraph mark/formatting), then
# remove two of the now-superfluous blank paragraphs that follow it.

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Test.`r") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Delete()

        # Remove two trailing empty paragraphs right after this one.
        $d.Paragraphs.Item($i + 1).Range.Delete()
        $d.Paragraphs.Item($i + 1).Range.Delete()

        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not find paragraph containing 'Test.'"
}
